# Actualización 10 de Mayo
#
# "Estadisticos 2P" (2nd partial) now has real grade statistics instead of
# the placeholder 0/blank figures, and "Estadisticos Final" is refreshed to
# match those new 2P numbers (its Blancos/Reprobados/Aprobados/Por_Apro and
# the resulting Promedio move together). "Estadisticos 1P" is left as-is.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param($ws, $row, $firstCol, [object[]]$values)
    $cols = @("A","B","C","D","E","F","G","H")
    $colIndex = [Array]::IndexOf($cols, $firstCol)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $colLetter = $cols[$colIndex + $i]
        $ws.Range($colLetter + $row).Value = $values[$i]
    }
}

# --- Estadisticos 2P ---------------------------------------------------
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")

Set-Row $ws2P 2 "D" @(0, 10, 29, 74.36, 7.2)
Set-Row $ws2P 3 "D" @(0, 4, 37, 90.24, 8.1)
Set-Row $ws2P 4 "D" @(0, 2, 23, 92, 7.4)
Set-Row $ws2P 5 "D" @(0, 6, 33, 84.62, 7.7)
Set-Row $ws2P 6 "D" @(0, 8, 27, 77.14, 7.7)
Set-Row $ws2P 7 "D" @(0, 9, 27, 75, 7.6)

# --- Estadisticos Final -------------------------------------------------
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")

Set-Row $wsFinal 2 "H" @(7.3)
Set-Row $wsFinal 3 "E" @(4, 37, 90.24, 8.2)
Set-Row $wsFinal 5 "E" @(6, 33, 84.62, 7.5)
Set-Row $wsFinal 6 "E" @(8, 27, 77.14, 7.5)
Set-Row $wsFinal 7 "E" @(9, 27, 75, 7.5)
